$wb = $excel.ActiveWorkbook

# Set value "A" in Sheet3!A1 (text, will become a shared string)
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws3 = $wb.Worksheets.Item("Sheet3")

$ws3.Range("A1").Value = "A"

# Update selections: Sheet1 and Sheet2 -> A1 selected (was A1:B2)
$ws1.Range("A1").Select() | Out-Null
$ws2.Range("A1").Select() | Out-Null

# Sheet3 -> A2 selected/active (was A1:B2)
$ws3.Range("A2").Select() | Out-Null
